{"js": "// Find the bullet that tells students to paste their code, and extend it\n// so it also mentions the option of attaching a .ino file, per the commit\n// message: \"students may attach .ino files to lab report\".\nconst body = context.document.body;\n\nconst results = body.search(\"copy and paste your code.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence 'copy and paste your code.'\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\n  \"copy and paste your code below. You may instead attach a .ino file if you prefer.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Extend the \"copy and paste your code.\" bullet so it also tells students\n# they may attach a .ino file instead, per the commit message:\n# \"students may attach .ino files to lab report\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"copy and paste your code.\"\n$find.Replacement.Text = \"copy and paste your code below. You may instead attach a .ino file if you prefer.\"\n\n# wdFindContinue = 1, wdReplaceOne = 1 (last arg)\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n"}
